$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New species columns
$ws.Range("F1").Value = "Alaria.marginata"
$ws.Range("G1").Value = "Costaria.costata"

# New observation rows (row 11: 2023-05-08, row 12: 2023-06-06)
$ws.Range("A11").Value = 2023
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = "y"
$ws.Range("E11").Value = "y"

$ws.Range("A12").Value = 2023
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = "y"
$ws.Range("E12").Value = "y"

$ws.Range("F12").Select()
